# "Modifiche pre trailer 2"
# Update the Instagram upload-plan sheet: add a "pubblicato" tracking column,
# turn a couple of the week-1/2 dates into real dates, and rewrite several
# "Tipo di Contenuto" entries to reflect what was actually posted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F: "pubblicato" (published?) -------------------------------
$ws.Range("F1").Value = "pubblicato"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("F1").VerticalAlignment = -4160     # xlTop
$ws.Range("F1").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("F1").Borders.Item(10).LineStyle = 1  # xlEdgeRight

$ws.Range("F2").Value = "no"
$ws.Range("F3").Value = "no"
$ws.Range("F4").Value = "no"
$ws.Range("F5").Value = "si"

# --- Real dates for two of the rows (were plain text before) ---------------
$ws.Range("B4").Value = [datetime]"2024-10-31"
$ws.Range("B5").Value = [datetime]"2024-11-02"

# --- Updated "Tipo di Contenuto" text in several rows -----------------------
$ws.Range("E6").Value  = "Reel con spiegazione"
$ws.Range("E8").Value  = "immagine"
$ws.Range("E9").Value  = "reel con spiegazione"
$ws.Range("E11").Value = "Reel di intrattenimento (es. tendenze o meme del settore)"
$ws.Range("E12").Value = "immagine"
$ws.Range("E13").Value = "immagine"

# E11 / E13 get an underline to flag them as re-used/duplicated entries.
$ws.Range("E11").Font.Underline = 2   # xlUnderlineStyleSingle
$ws.Range("E13").Font.Underline = 2   # xlUnderlineStyleSingle

# --- Page setup tweaks -------------------------------------------------------
$ws.PageSetup.PaperSize = 9    # xlPaperA4
$ws.PageSetup.Orientation = 1  # xlPortrait

Write-Output "done"
